# Applies the Lab4.docx edit:
#   1. merges the "Горностай Богдан " / ", КН-108" runs into one run;
#   2. relocates the hidden "_GoBack" bookmark from the blank paragraph
#      before "Опис програми" to the start of the "Програма створює..."
#      paragraph;
#   3. rewords that paragraph (drops the "аудіотек(у/и)" phrasing for
#      "каталог розкладів" / "маршрути");
#   4. updates that paragraph's indent from right/firstLine to left/right.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Горностай Богдан " + ", КН-108" -> single run "Горностай Богдан , КН-108"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Горностай Богдан , КН-108", $false, $false, $false, $false, $false,
    $true, 1, $false, "Горностай Богдан , КН-108", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Move the hidden "_GoBack" bookmark from the empty paragraph (right
#    before "Опис програми") down onto the start of the "Програма створює"
#    paragraph. Word keeps a single "_GoBack" bookmark, so re-adding it at
#    the new location removes it from the old one automatically.
# ---------------------------------------------------------------------------
$goBackPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $r = $p.Range.Duplicate
    $hit = $r.Find.Execute(
        "Програма створює", $false, $false, $false, $false, $false,
        $true, 1, $false, "", 0)
    if ($hit) {
        $goBackPara = $p
        break
    }
}

$target = $goBackPara.Range.Duplicate
$target.Collapse(1)
$d.Bookmarks.Add("_GoBack", $target) | Out-Null

# ---------------------------------------------------------------------------
# 3) Rewrite the paragraph text: drop "аудіотек(у|и)" wording in favour of
#    the "маршрут"/"каталог розкладів" wording.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Програма створює вашу власну аудіотеку за допомогою класу",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Програма створює вашу власний каталог розкладів за допомогою класу",
    2) | Out-Null

$d.Content.Find.Execute(
    " ваші аудіотеки у файл вибраний користувачем",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " ваші маршрути у файл вибраний користувачем",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Paragraph indent: right=692/firstLine=706 -> left=1038/right=692
# ---------------------------------------------------------------------------
$goBackPara.Format.LeftIndent = 51.9
$goBackPara.Format.FirstLineIndent = 0
